$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the cells we are about to update to remain TEXT cells (matching the
# original t="inlineStr" string cells), rather than letting Excel auto-convert
# numeric-looking strings (prices) or percent-looking strings into numbers.
$textCells = @(
    "D2", "E2", "G2", "D3", "E3", "G3", "D4", "E4", "G4", "D5",
    "E5", "G5", "D6", "E6", "G6", "D7", "E7", "G7", "D8", "E8",
    "G8", "E9", "G9", "D10", "E10", "G10", "D11", "E11", "G11", "D12",
    "E12", "G12", "D13", "E13", "G13", "D14", "E14", "G14", "E15", "G15",
    "D16", "E16", "G16", "D17", "E17", "G17", "D18", "E18", "G18", "E19",
    "G19", "E20", "G20", "D21", "E21", "G21", "D22", "E22", "G22", "D23",
    "E23", "G23", "D24", "E24", "G24", "D25", "E25", "G25", "D26", "E26",
    "G26", "G27", "G28", "G29", "G30", "G31", "G32", "G33", "G34", "G35",
    "G36", "G37", "G38", "D39", "E39", "G39", "D40", "E40", "G40", "D41",
    "E41", "G41", "D42", "E42", "G42", "D43", "E43", "G43", "D44", "E44",
    "G44", "D45", "E45", "G45", "D46", "E46", "G46", "E47", "G47", "D48",
    "E48", "G48", "E49", "G49", "E50", "G50", "E51", "G51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values from the crypto price/volume/hour refresh.
$ws.Range("D2").Value = '331.33'
$ws.Range("E2").Value = '0.11%'
$ws.Range("G2").Value = '16'

$ws.Range("D3").Value = '41.58'
$ws.Range("E3").Value = '0.86%'
$ws.Range("G3").Value = '16'

$ws.Range("D4").Value = '5.694'
$ws.Range("E4").Value = '-0.92%'
$ws.Range("G4").Value = '16'

$ws.Range("D5").Value = '0.08354'
$ws.Range("E5").Value = '2.99%'
$ws.Range("G5").Value = '16'

$ws.Range("D6").Value = '8.790'
$ws.Range("E6").Value = '1.28%'
$ws.Range("G6").Value = '16'

$ws.Range("D7").Value = '2.009'
$ws.Range("E7").Value = '1.51%'
$ws.Range("G7").Value = '16'

$ws.Range("D8").Value = '4.518'
$ws.Range("E8").Value = '0.27%'
$ws.Range("G8").Value = '16'

$ws.Range("E9").Value = '-1.24%'
$ws.Range("G9").Value = '16'

$ws.Range("D10").Value = '0.9271'
$ws.Range("E10").Value = '-0.07%'
$ws.Range("G10").Value = '16'

$ws.Range("D11").Value = '0.1294'
$ws.Range("E11").Value = '1.11%'
$ws.Range("G11").Value = '16'

$ws.Range("D12").Value = '0.1973'
$ws.Range("E12").Value = '0.92%'
$ws.Range("G12").Value = '16'

$ws.Range("D13").Value = '0.09570'
$ws.Range("E13").Value = '4.18%'
$ws.Range("G13").Value = '16'

$ws.Range("D14").Value = '0.03889'
$ws.Range("E14").Value = '4.62%'
$ws.Range("G14").Value = '16'

$ws.Range("E15").Value = '0.96%'
$ws.Range("G15").Value = '16'

$ws.Range("D16").Value = '0.001309'
$ws.Range("E16").Value = '1.22%'
$ws.Range("G16").Value = '16'

$ws.Range("D17").Value = '0.006091'
$ws.Range("E17").Value = '-2.24%'
$ws.Range("G17").Value = '16'

$ws.Range("D18").Value = '3.438'
$ws.Range("E18").Value = '2.00%'
$ws.Range("G18").Value = '16'

$ws.Range("E19").Value = '1.18%'
$ws.Range("G19").Value = '16'

$ws.Range("E20").Value = '-7.47%'
$ws.Range("G20").Value = '16'

$ws.Range("D21").Value = '0.1364'
$ws.Range("E21").Value = '-0.93%'
$ws.Range("G21").Value = '16'

$ws.Range("D22").Value = '0.2511'
$ws.Range("E22").Value = '-3.79%'
$ws.Range("G22").Value = '16'

$ws.Range("D23").Value = '0.04436'
$ws.Range("E23").Value = '0.24%'
$ws.Range("G23").Value = '16'

$ws.Range("D24").Value = '0.001278'
$ws.Range("E24").Value = '1.67%'
$ws.Range("G24").Value = '16'

$ws.Range("D25").Value = '0.004412'
$ws.Range("E25").Value = '0.47%'
$ws.Range("G25").Value = '16'

$ws.Range("D26").Value = '0.0001202'
$ws.Range("E26").Value = '-3.08%'
$ws.Range("G26").Value = '16'

$ws.Range("G27").Value = '16'

$ws.Range("G28").Value = '16'

$ws.Range("G29").Value = '16'

$ws.Range("G30").Value = '16'

$ws.Range("G31").Value = '16'

$ws.Range("G32").Value = '16'

$ws.Range("G33").Value = '16'

$ws.Range("G34").Value = '16'

$ws.Range("G35").Value = '16'

$ws.Range("G36").Value = '16'

$ws.Range("G37").Value = '16'

$ws.Range("G38").Value = '16'

$ws.Range("D39").Value = '0.02812'
$ws.Range("E39").Value = '-0.40%'
$ws.Range("G39").Value = '16'

$ws.Range("D40").Value = '0.05561'
$ws.Range("E40").Value = '0.42%'
$ws.Range("G40").Value = '16'

$ws.Range("D41").Value = '0.007804'
$ws.Range("E41").Value = '1.46%'
$ws.Range("G41").Value = '16'

$ws.Range("D42").Value = '0.1437'
$ws.Range("E42").Value = '1.27%'
$ws.Range("G42").Value = '16'

$ws.Range("D43").Value = '0.009311'
$ws.Range("E43").Value = '-5.51%'
$ws.Range("G43").Value = '16'

$ws.Range("D44").Value = '0.002161'
$ws.Range("E44").Value = '-2.83%'
$ws.Range("G44").Value = '16'

$ws.Range("D45").Value = '0.01103'
$ws.Range("E45").Value = '-7.13%'
$ws.Range("G45").Value = '16'

$ws.Range("D46").Value = '0.00007006'
$ws.Range("E46").Value = '3.16%'
$ws.Range("G46").Value = '16'

$ws.Range("E47").Value = '0.17%'
$ws.Range("G47").Value = '16'

$ws.Range("D48").Value = '0.003528'
$ws.Range("E48").Value = '14.84%'
$ws.Range("G48").Value = '16'

$ws.Range("E49").Value = '0.11%'
$ws.Range("G49").Value = '16'

$ws.Range("E50").Value = '0.17%'
$ws.Range("G50").Value = '16'

$ws.Range("E51").Value = '0.17%'
$ws.Range("G51").Value = '16'

